$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# New dates (column A) for rows 692-705
$dates = @(
    "2022-02-26",
    "2022-02-27",
    "2022-02-28",
    "2022-03-01",
    "2022-03-02",
    "2022-03-03",
    "2022-03-04",
    "2022-03-05",
    "2022-03-06",
    "2022-03-07",
    "2022-03-08",
    "2022-03-09",
    "2022-03-10",
    "2022-03-11"
)

# New Scheduled flights (column B) and Tracked flights (column C) values
$scheduled = @(54, 64, 72, 65, 61, 72, 76, 56, 70, 61, 62, 66, 76, 80)
$tracked   = @(51, 60, 70, 65, 59, 68, 71, 56, 66, 58, 61, 62, 73, 73)

$lastRow = 691
$startRow = 692

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i

    # Clone formatting from the last existing data row before writing values,
    # so new cells keep the same style indices (text/number/percent) instead
    # of Excel auto-detecting a date/general format for the fresh cells.
    $ws.Cells.Item($lastRow, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)
    $ws.Cells.Item($row, 1).Value = $dates[$i]

    $ws.Cells.Item($lastRow, 2).Copy()
    $ws.Cells.Item($row, 2).PasteSpecial(-4122)
    $ws.Cells.Item($row, 2).Value = $scheduled[$i]

    $ws.Cells.Item($lastRow, 3).Copy()
    $ws.Cells.Item($row, 3).PasteSpecial(-4122)
    $ws.Cells.Item($row, 3).Value = $tracked[$i]

    $ws.Cells.Item($lastRow, 4).Copy()
    $ws.Cells.Item($row, 4).PasteSpecial(-4122)
    $ws.Cells.Item($row, 4).Formula = "=C$row/B$row"
}

$excel.CutCopyMode = $false

# Update the visible viewport / active selection to match the post-edit state
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 320
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H701").Select()
